$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, copying the formatting from G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Rows where Save = 1 (based on the original diff data)
$saveOneRows = @(5, 8, 22)

for ($r = 2; $r -le 61; $r++) {
    if ($saveOneRows -contains $r) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
